$wb = $excel.ActiveWorkbook

# --- Section_A ---
$ws1 = $wb.Worksheets.Item("Section_A")

# Extend the time-slot grid: copy style of existing time-label cell A7 into new A8:A12
for ($r = 8; $r -le 12; $r++) {
    $ws1.Range("A7").Copy($ws1.Cells.Item($r, 1))
}

# Updated values for existing cells
$ws1.Range("B2").Value = "CS304"
$ws1.Range("C2").Value = "CS304"
$ws1.Range("F2").Value = "CS303"
$ws1.Range("F3").Value = "Free"
$ws1.Range("C5").Value = "Free"
$ws1.Range("D5").Value = "CS309"
$ws1.Range("F5").Value = "Free"
$ws1.Range("B6").Value = "CS303"
$ws1.Range("C6").Value = "Free"
$ws1.Range("D6").Value = "Free"
$ws1.Range("B7").Value = "CS309"
$ws1.Range("D7").Value = "CS303"
$ws1.Range("E7").Value = "CS304"
$ws1.Range("F7").Value = "CS309"

# New cells (rows 8-12)
$ws1.Range("A8").Value = "12:00-13:00"
$ws1.Range("B8").Value = "Free"
$ws1.Range("C8").Value = "Free"
$ws1.Range("D8").Value = "CS304 (Tutorial)"
$ws1.Range("E8").Value = "Free"
$ws1.Range("F8").Value = "Free"
$ws1.Range("A9").Value = "13:00-14:00"
$ws1.Range("B9").Value = "Free"
$ws1.Range("C9").Value = "Free"
$ws1.Range("D9").Value = "CS303 (Tutorial)"
$ws1.Range("E9").Value = "Free"
$ws1.Range("F9").Value = "Free"
$ws1.Range("A10").Value = "15:30-16:30"
$ws1.Range("B10").Value = "Free"
$ws1.Range("C10").Value = "Free"
$ws1.Range("D10").Value = "Free"
$ws1.Range("E10").Value = "Free"
$ws1.Range("F10").Value = "CS309 (Tutorial)"
$ws1.Range("A11").Value = "16:30-17:30"
$ws1.Range("B11").Value = "Free"
$ws1.Range("C11").Value = "Free"
$ws1.Range("D11").Value = "Free"
$ws1.Range("E11").Value = "Free"
$ws1.Range("F11").Value = "Free"
$ws1.Range("A12").Value = "17:30-18:30"
$ws1.Range("B12").Value = "Free"
$ws1.Range("C12").Value = "Free"
$ws1.Range("D12").Value = "Free"
$ws1.Range("E12").Value = "Free"
$ws1.Range("F12").Value = "Free"

# --- Section_B ---
$ws2 = $wb.Worksheets.Item("Section_B")

# Extend the time-slot grid: copy style of existing time-label cell A7 into new A8:A12
for ($r = 8; $r -le 12; $r++) {
    $ws2.Range("A7").Copy($ws2.Cells.Item($r, 1))
}

# Updated values for existing cells
$ws2.Range("C2").Value = "CS304"
$ws2.Range("D2").Value = "Free"
$ws2.Range("E2").Value = "Free"
$ws2.Range("B3").Value = "CS303"
$ws2.Range("D3").Value = "Free"
$ws2.Range("F3").Value = "Free"
$ws2.Range("B5").Value = "CS309"
$ws2.Range("C5").Value = "CS303"
$ws2.Range("D5").Value = "Free"
$ws2.Range("E5").Value = "Free"
$ws2.Range("F5").Value = "CS309"
$ws2.Range("B6").Value = "CS304"
$ws2.Range("C6").Value = "CS309"
$ws2.Range("F7").Value = "Free"

# New cells (rows 8-12)
$ws2.Range("A8").Value = "12:00-13:00"
$ws2.Range("B8").Value = "Free"
$ws2.Range("C8").Value = "Free"
$ws2.Range("D8").Value = "Free"
$ws2.Range("E8").Value = "Free"
$ws2.Range("F8").Value = "Free"
$ws2.Range("A9").Value = "13:00-14:00"
$ws2.Range("B9").Value = "Free"
$ws2.Range("C9").Value = "Free"
$ws2.Range("D9").Value = "Free"
$ws2.Range("E9").Value = "CS309 (Tutorial)"
$ws2.Range("F9").Value = "Free"
$ws2.Range("A10").Value = "15:30-16:30"
$ws2.Range("B10").Value = "Free"
$ws2.Range("C10").Value = "Free"
$ws2.Range("D10").Value = "Free"
$ws2.Range("E10").Value = "CS303 (Tutorial)"
$ws2.Range("F10").Value = "Free"
$ws2.Range("A11").Value = "16:30-17:30"
$ws2.Range("B11").Value = "Free"
$ws2.Range("C11").Value = "Free"
$ws2.Range("D11").Value = "Free"
$ws2.Range("E11").Value = "Free"
$ws2.Range("F11").Value = "Free"
$ws2.Range("A12").Value = "17:30-18:30"
$ws2.Range("B12").Value = "Free"
$ws2.Range("C12").Value = "CS304 (Tutorial)"
$ws2.Range("D12").Value = "Free"
$ws2.Range("E12").Value = "Free"
$ws2.Range("F12").Value = "Free"

# --- Course_Summary ---
$ws3 = $wb.Worksheets.Item("Course_Summary")

# New header cells G1/H1 need the bold/centered header style; copy it from F1 first
$ws3.Range("F1").Copy($ws3.Range("G1"))
$ws3.Range("F1").Copy($ws3.Range("H1"))

# Updated values for existing cells
$ws3.Range("E1").Value = "Lectures/Week"
$ws3.Range("F1").Value = "Tutorials/Week"
$ws3.Range("E2").Value = 3
$ws3.Range("F2").Value = 1
$ws3.Range("E3").Value = 3
$ws3.Range("F3").Value = 1
$ws3.Range("E4").Value = 3
$ws3.Range("F4").Value = 1
$ws3.Range("F5").Value = 0

# New cells (columns G,H)
$ws3.Range("G1").Value = "Total Credits"
$ws3.Range("H1").Value = "Instructor"
$ws3.Range("G2").Value = 4
$ws3.Range("H2").Value = "Dr. Sunil C K"
$ws3.Range("G3").Value = 6
$ws3.Range("H3").Value = "Dr. Animesh Roy"
$ws3.Range("G4").Value = 4
$ws3.Range("H4").Value = "Dr. Krishnendu"
$ws3.Range("G5").Value = 4
$ws3.Range("H5").Value = "Dr. Pramod"
